$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D19,D20,D22,D23,D24,D25: switch border/format from the "no-right-border"
# style (used elsewhere in this block) to the full-box-border style that D21
# already uses (copy formats only, values are untouched). ---
$ws.Range("D21").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D25").PasteSpecial(-4122)

# --- D26:D43 (probateFormsRW01.feature "SmokeTest" column): No -> Yes ---
$ws.Range("D26").Value = "Yes"
$ws.Range("D27").Value = "Yes"
$ws.Range("D28").Value = "Yes"
$ws.Range("D29").Value = "Yes"
$ws.Range("D30").Value = "Yes"
$ws.Range("D31").Value = "Yes"
$ws.Range("D32").Value = "Yes"
$ws.Range("D33").Value = "Yes"
$ws.Range("D34").Value = "Yes"
$ws.Range("D35").Value = "Yes"
$ws.Range("D36").Value = "Yes"
$ws.Range("D37").Value = "Yes"
$ws.Range("D38").Value = "Yes"
$ws.Range("D39").Value = "Yes"
$ws.Range("D40").Value = "Yes"
$ws.Range("D41").Value = "Yes"
$ws.Range("D42").Value = "Yes"
$ws.Range("D43").Value = "Yes"

# --- D74:D84: switch from full-box-border style back to the no-right-border
# style that D44 (and the rest of the RW02/RW03 block) already uses. ---
$ws.Range("D44").Copy()
$ws.Range("D74").PasteSpecial(-4122)
$ws.Range("D75").PasteSpecial(-4122)
$ws.Range("D76").PasteSpecial(-4122)
$ws.Range("D77").PasteSpecial(-4122)
$ws.Range("D78").PasteSpecial(-4122)
$ws.Range("D79").PasteSpecial(-4122)
$ws.Range("D80").PasteSpecial(-4122)
$ws.Range("D81").PasteSpecial(-4122)
$ws.Range("D82").PasteSpecial(-4122)
$ws.Range("D83").PasteSpecial(-4122)
$ws.Range("D84").PasteSpecial(-4122)

# --- C96:C106 (probateFormsRW05.feature "ScenarioName" column): de-duplicate
# its cell format onto the equivalent xf already used by C85 ---
$ws.Range("C85").Copy()
$ws.Range("C96").PasteSpecial(-4122)
$ws.Range("C97").PasteSpecial(-4122)
$ws.Range("C98").PasteSpecial(-4122)
$ws.Range("C99").PasteSpecial(-4122)
$ws.Range("C100").PasteSpecial(-4122)
$ws.Range("C101").PasteSpecial(-4122)
$ws.Range("C102").PasteSpecial(-4122)
$ws.Range("C103").PasteSpecial(-4122)
$ws.Range("C104").PasteSpecial(-4122)
$ws.Range("C105").PasteSpecial(-4122)
$ws.Range("C106").PasteSpecial(-4122)

# --- D96:D106 (probateFormsRW05.feature "SmokeTest" column): Yes -> No ---
$ws.Range("D96").Value = "No"
$ws.Range("D97").Value = "No"
$ws.Range("D98").Value = "No"
$ws.Range("D99").Value = "No"
$ws.Range("D100").Value = "No"
$ws.Range("D101").Value = "No"
$ws.Range("D102").Value = "No"
$ws.Range("D103").Value = "No"
$ws.Range("D104").Value = "No"
$ws.Range("D105").Value = "No"
$ws.Range("D106").Value = "No"

# --- Move the active selection to C26 (and let the view scroll naturally,
# dropping the stale topLeftCell="A92" pin) ---
$ws.Range("C26").Select()
